# Generate Report for Handback
# Row 17 in both the "zh-cn" and "de-de" sheets previously re-used the same
# Correspond Handoff Datetime / Correspond Handback DateTime values as row 18
# (a report-generation bug). Regenerate fresh, unique timestamps for row 17
# on both sheets, leaving row 18 (and everything else) untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D17").Value = "2016-03-10 05:55:14"
$wsZhCn.Range("G17").Value = "2016-03-10 05:55:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D17").Value = "2016-03-10 05:55:22"
$wsDeDe.Range("G17").Value = "2016-03-10 05:56:11"
